# Update the SexRatio (column N) values for rows 4-17 to a constant 0.64,
# overwriting the existing formulas in N12 and N16 with a plain static value,
# as per reviewer comments.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MinCount_summary_KZ-withimm")

$ws.Range("N4").Value = 0.64
$ws.Range("N5").Value = 0.64
$ws.Range("N6").Value = 0.64
$ws.Range("N7").Value = 0.64
$ws.Range("N8").Value = 0.64
$ws.Range("N9").Value = 0.64
$ws.Range("N10").Value = 0.64
$ws.Range("N11").Value = 0.64
$ws.Range("N12").Value = 0.64
$ws.Range("N13").Value = 0.64
$ws.Range("N14").Value = 0.64
$ws.Range("N15").Value = 0.64
$ws.Range("N16").Value = 0.64
$ws.Range("N17").Value = 0.64

# Match the resulting view state: scrolled/selected cell moved from W34 to Q34.
$ws.Activate()
$ws.Range("Q34").Select()
